$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 92
$ws.Range("I5").Value = 92
$ws.Range("K5").Value = 92
$ws.Range("M5").Value = 23
$ws.Range("H21").Value = 20008.5
$ws.Range("I21").Value = 20008.5
$ws.Range("K21").Value = 20008.5
$ws.Range("M21").Value = -19540.5
$ws.Range("H23").Value = 20008.5
$ws.Range("I23").Value = 20008.5
$ws.Range("K23").Value = 20008.5
$ws.Range("M23").Value = -19774.5
$ws.Range("H29").Value = 1500
$ws.Range("J29").Value = 1500
$ws.Range("L29").Value = 4500
$ws.Range("N29").Value = -5062
$ws.Range("H38").Value = 1818.1428
$ws.Range("I38").Value = 1818.1428
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 5454.428400000001
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -5082.428400000001
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 728.63635
$ws.Range("J58").Value = 675
$ws.Range("L58").Value = 2025
$ws.Range("N58").Value = -2325
$ws.Range("H69").Value = 2500
$ws.Range("I69").Value = 2500
$ws.Range("K69").Value = 7500
$ws.Range("M69").Value = -6626
$ws.Range("H72").Value = 2500
$ws.Range("I72").Value = 2500
$ws.Range("K72").Value = 22500
$ws.Range("M72").Value = -18132
$ws.Range("H88").Value = 2444.3333
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 2444.3333
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H93").Value = 35000
$ws.Range("J93").Value = 35000
$ws.Range("L93").Value = 35000
$ws.Range("N93").Value = -39992
$ws.Range("H99").Value = 1700
$ws.Range("J99").Value = 1700
$ws.Range("L99").Value = 5100
$ws.Range("N99").Value = -8096
$ws.Range("H115").Value = 400
$ws.Range("I115").Value = 400
$ws.Range("K115").Value = 1200
$ws.Range("M115").Value = 367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1836.8334
$ws.Range("I2").Value = 1952.5
$ws.Range("J2").Value = 1605.5
$ws.Range("K2").Value = 1952.5
$ws.Range("L2").Value = 1605.5
$ws.Range("M2").Value = -1839.5
$ws.Range("N2").Value = -1831.5
$ws.Range("H5").Value = 225.66667
$ws.Range("I5").Value = 237.5
$ws.Range("K5").Value = 237.5
$ws.Range("M5").Value = -125.5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H63").Value = 1550
$ws.Range("I63").Value = 1550
$ws.Range("K63").Value = 1550
$ws.Range("M63").Value = -864
$ws.Range("H66").Value = 1550
$ws.Range("I66").Value = 1550
$ws.Range("K66").Value = 7750
$ws.Range("M66").Value = -4318
$ws.Range("H96").Value = 29629.125
$ws.Range("J96").Value = 29629.125
$ws.Range("L96").Value = 29629.125
$ws.Range("N96").Value = -35121.125
$ws.Range("H102").Value = 2133.2856
$ws.Range("I102").Value = 2133.2856
$ws.Range("K102").Value = 2133.2856
$ws.Range("M102").Value = -511.2856000000002
$ws.Range("H116").Value = 1836.8334
$ws.Range("I116").Value = 1952.5
$ws.Range("J116").Value = 1605.5
$ws.Range("K116").Value = 1952.5
$ws.Range("L116").Value = 1605.5
$ws.Range("M116").Value = 341.5
$ws.Range("N116").Value = -6193.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1836.8334
$ws.Range("I3").Value = 1952.5
$ws.Range("J3").Value = 1605.5
$ws.Range("K3").Value = 1952.5
$ws.Range("L3").Value = 1605.5
$ws.Range("M3").Value = -1838.5
$ws.Range("N3").Value = -1833.5
$ws.Range("H4").Value = 225.66667
$ws.Range("I4").Value = 237.5
$ws.Range("K4").Value = 237.5
$ws.Range("M4").Value = -122.5
$ws.Range("H134").Value = 3001.5
$ws.Range("I134").Value = 2483.125
$ws.Range("K134").Value = 7449.375
$ws.Range("M134").Value = -4914.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 518.2
$ws.Range("I22").Value = 580.3333
$ws.Range("K22").Value = 580.3333
$ws.Range("M22").Value = -230.3333
$ws.Range("H31").Value = 1071.091
$ws.Range("I31").Value = 1071.091
$ws.Range("K31").Value = 1071.091
$ws.Range("M31").Value = -776.0909999999999
$ws.Range("H34").Value = 1071.091
$ws.Range("I34").Value = 1071.091
$ws.Range("K34").Value = 1071.091
$ws.Range("M34").Value = -869.0909999999999
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H96").Value = 12820.5
$ws.Range("J96").Value = 12820.5
$ws.Range("L96").Value = 12820.5
$ws.Range("N96").Value = -18312.5
$ws.Range("H105").Value = 6201
$ws.Range("I105").Value = 6201
$ws.Range("K105").Value = 6201
$ws.Range("M105").Value = -4454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 100
$ws.Range("I50").Value = 100
$ws.Range("K50").Value = 300
$ws.Range("M50").Value = 181
$ws.Range("H53").Value = 100
$ws.Range("I53").Value = 100
$ws.Range("K53").Value = 300
$ws.Range("M53").Value = 181
$ws.Range("H80").Value = 4000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 6000
$ws.Range("L80").Value = 18000
$ws.Range("M80").Value = -5064
$ws.Range("N80").Value = -19872
$ws.Range("H83").Value = 4000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 18000
$ws.Range("L83").Value = 54000
$ws.Range("M83").Value = -13320
$ws.Range("N83").Value = -63360
$ws.Range("H131").Value = 3665.3845
$ws.Range("J131").Value = 3800
$ws.Range("L131").Value = 11400
$ws.Range("N131").Value = -21480
$ws.Range("H140").Value = 995.75
$ws.Range("I140").Value = 995.75
$ws.Range("K140").Value = 2987.25
$ws.Range("M140").Value = 2192.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 2686.6667
$ws.Range("I20").Value = 21
$ws.Range("J20").Value = 40006
$ws.Range("K20").Value = 21
$ws.Range("L20").Value = 40006
$ws.Range("M20").Value = 224
$ws.Range("N20").Value = -40496
$ws.Range("H80").Value = 1924.8334
$ws.Range("I80").Value = 1812.25
$ws.Range("J80").Value = 2150
$ws.Range("K80").Value = 1812.25
$ws.Range("L80").Value = 2150
$ws.Range("M80").Value = -814.25
$ws.Range("N80").Value = -4146
$ws.Range("H83").Value = 1924.8334
$ws.Range("I83").Value = 1812.25
$ws.Range("J83").Value = 2150
$ws.Range("K83").Value = 9061.25
$ws.Range("L83").Value = 10750
$ws.Range("M83").Value = -4069.25
$ws.Range("N83").Value = -20734
$ws.Range("H95").Value = 14750
$ws.Range("J95").Value = 14750
$ws.Range("L95").Value = 14750
$ws.Range("N95").Value = -20242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1150
$ws.Range("I46").Value = 1150
$ws.Range("K46").Value = 1150
$ws.Range("M46").Value = -962
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 23103.75
$ws.Range("I41").Value = 28000
$ws.Range("J41").Value = 18207.5
$ws.Range("K41").Value = 28000
$ws.Range("L41").Value = 18207.5
$ws.Range("M41").Value = -27610
$ws.Range("N41").Value = -18987.5
$ws.Range("J81").Value = 500
$ws.Range("L81").Value = 1000
$ws.Range("N81").Value = -3122
$ws.Range("J84").Value = 500
$ws.Range("L84").Value = 5000
$ws.Range("N84").Value = -15608
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H136").Value = 92539.09
$ws.Range("I136").Value = 760.125
$ws.Range("K136").Value = 2280.375
$ws.Range("M136").Value = 269.625

Write-Host "All changes applied."